$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.143.32"
$ws.Range("E2").Value = "  -2.97%  "

$ws.Range("D3").Value = "1.596.59"
$ws.Range("E3").Value = "  -3.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3754"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3654"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.004"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.275"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08060"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.624"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.553"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001262"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.63%  "

$ws.Range("E17").Value = "  -3.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06809"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.594"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.56%  "

$ws.Range("D24").Value = "23.153.17"
$ws.Range("E24").Value = "  -2.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.359"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.925"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.238"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.447"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.185"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.20%  "

$ws.Range("D33").Value = "1.771.40"
$ws.Range("E33").Value = "  -3.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9728"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07743"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02777"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.279"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2536"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08840"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.392"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7150"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6615"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.03%  "

$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.306"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.961"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07980"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.167"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.19%  "
